$d = $word.ActiveDocument

$pairs = @(
    @{old="18×43="; new="89×44="},
    @{old="55×99="; new="21×26="},
    @{old="71×25="; new="46×55="},
    @{old="64×37="; new="43×80="},
    @{old="28×83="; new="62×18="},
    @{old="17×52="; new="19×32="},
    @{old="22×27="; new="26×99="},
    @{old="66×91="; new="58×66="},
    @{old="33×66="; new="77×57="},
    @{old="39×48="; new="68×89="},
    @{old="41×53="; new="50×16="},
    @{old="93×85="; new="35×67="},
    @{old="84×73="; new="61×12="},
    @{old="80×50="; new="64×81="},
    @{old="39×96="; new="74×49="},
    @{old="54×84="; new="76×17="},
    @{old="50×19="; new="56×66="},
    @{old="85×25="; new="77×61="},
    @{old="69×71="; new="32×17="},
    @{old="34×41="; new="70×51="},
    @{old="37×26="; new="88×48="},
    @{old="85×92="; new="94×44="},
    @{old="20×49="; new="98×94="},
    @{old="93×32="; new="50×14="},
    @{old="19×67="; new="15×49="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
